# Update Sheet1 with new bioassay data (new date + new measurements).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New raw data rows (A2:H6) replacing the old ones.
$data = @(
    @(45113, 4.8896115217396465, 3.9476581455361317, 23.639428428834751, 6.5241821534878355, 4.186488657491406,  14.757457401238165, 28.500277760066652),
    @(45113, 4.9431171487971159, 4.1271753721004609, 25.191950096368846, 5.2087956578165286, 4.0054451095318475, 23.010188106063957, 28.024101383826551),
    @(45113, 4.6960663481837503, 3.8811016984526123, 25.470488204324724, 4.6980278967165479, 4.2109460521046769, 24.539647899358865, 28.965265419875561),
    @(45113, 5.0553619656018114, 4.2327127993789944, 28.028458608692059, 4.4450988467522272, 3.9548254964979823, 26.940918788531086, 28.604282903445533),
    @(45113, 4.7132711724114937, 4.5629542346052281, 30.397299741116068, 5.3310334382232423, 4.0650677892302456, 29.255218623775907, 23.086322575178144)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($col = 1; $col -le 8; $col++) {
        $ws1.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}

# Sheet1 rows 7-9 (avg / change-from-t0 / change-from-control) are formulas
# and recalc automatically. Sheet2 stores a plain copy of those results, so
# push the refreshed numbers over there too.
$ws2.Range("B2").Value = $ws1.Range("B7").Value2
$ws2.Range("B3").Value = $ws1.Range("C7").Value2
$ws2.Range("B4").Value = $ws1.Range("D7").Value2
$ws2.Range("B5").Value = $ws1.Range("E7").Value2
$ws2.Range("B6").Value = $ws1.Range("F7").Value2
$ws2.Range("B7").Value = $ws1.Range("G7").Value2
$ws2.Range("B8").Value = $ws1.Range("H7").Value2

$ws2.Range("C3").Value = $ws1.Range("C8").Value2
$ws2.Range("C4").Value = $ws1.Range("D8").Value2
$ws2.Range("C5").Value = $ws1.Range("E8").Value2
$ws2.Range("C6").Value = $ws1.Range("F8").Value2
$ws2.Range("C7").Value = $ws1.Range("G8").Value2
$ws2.Range("C8").Value = $ws1.Range("H8").Value2

$ws2.Range("D4").Value = $ws1.Range("D9").Value2
$ws2.Range("D5").Value = $ws1.Range("E9").Value2
$ws2.Range("D6").Value = $ws1.Range("F9").Value2
$ws2.Range("D7").Value = $ws1.Range("G9").Value2
$ws2.Range("D8").Value = $ws1.Range("H9").Value2

# Sheet2 becomes the active/selected sheet with D2:D8 selected (this also
# drops Sheet1's now-stale tabSelected/topLeftCell view attributes).
$ws2.Activate()
$ws2.Range("D2:D8").Select()
